$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Datensatz BTW 2025")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Kurzform"

$ws4.Range("A1").Value = $ws3.Range("D1").Value()
$ws4.Range("B1").Value = $ws3.Range("E1").Value()
$ws4.Range("C1").Value = $ws3.Range("F1").Value()
$ws4.Range("D1").Value = "These: Kurzform"

Write-Output $ws4.Range("A1").Value()
Write-Output $ws4.Range("B1").Value()
Write-Output $ws4.Range("C1").Value()
Write-Output $ws4.Range("D1").Value()
